# Apply the commit's workbook-state changes:
#  1. "Words" sheet: refresh the analysis_started_at timestamp on row 3 (H3)
#     to the latest recalculated value.
#  2. "Progress" sheet: append three new "ਗੁਬਾਰੀ" rows (11-13) that were
#     queued/selected for analysis, extending the used range from
#     A1:K10 to A1:K13.

$wb = $excel.ActiveWorkbook

# --- 1. Words sheet: nudge H3's cached timestamp -------------------------
$wordsWs = $wb.Worksheets.Item("Words")
$wordsWs.Range("H3").Value = 45914.2781740625

# --- 2. Progress sheet: append the three new rows -------------------------
$progressWs = $wb.Worksheets.Item("Progress")

$newRows = @(
    @{ Row = 11; Word = "ਗੁਬਾਰੀ"; Verse = "ਮਨਮੁਖ ਦੁਬਿਧਾ ਦੁਰਮਤਿ ਬਿਆਪੇ ਜਿਨ ਅੰਤਰਿ ਮੋਹ ਗੁਬਾਰੀ ॥"; Page = 507 },
    @{ Row = 12; Word = "ਗੁਬਾਰੀ"; Verse = "ਜਿਨ੍ਹ੍ਹਿ ਕੀਏ ਤਿਸਹਿ ਨ ਜਾਣਨੀ ਮਨਮੁਖਿ ਗੁਬਾਰੀ ॥"; Page = 788 },
    @{ Row = 13; Word = "ਗੁਬਾਰੀ"; Verse = "ਬਾਹਰਿ ਭਸਮ ਲੇਪਨ ਕਰੇ ਅੰਤਰਿ ਗੁਬਾਰੀ ॥"; Page = 1243 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $progressWs.Cells.Item($row, 1).Value = $r.Word      # A: word
    $progressWs.Cells.Item($row, 2).Value = $r.Word      # B: word_key_norm
    $progressWs.Cells.Item($row, 3).Value = ""           # C: word_index (blank)
    $progressWs.Cells.Item($row, 4).Value = $r.Verse     # D: verse
    $progressWs.Cells.Item($row, 5).Value = $r.Page      # E: page_number
    $progressWs.Cells.Item($row, 6).Value = $true        # F: selected_for_analysis
    $progressWs.Cells.Item($row, 7).Value = 45914.29202659722  # G: selected_at
    $progressWs.Cells.Item($row, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $progressWs.Cells.Item($row, 8).Value = "not started" # H: status
    $progressWs.Cells.Item($row, 9).Value = ""           # I: completed_at (blank)
    $progressWs.Cells.Item($row, 10).Value = ""          # J: reanalyzed_count (blank)
    $progressWs.Cells.Item($row, 11).Value = ""          # K: last_reanalyzed_at (blank)
}
